$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper scratch cell (well outside the used range) used to stamp text
# values via Copy / PasteSpecial(values) so the literal strings land as
# plain shared-string cells (t="s") without Excel's "smart" number/date
# auto-conversion forcing a new number-format style onto the cell.
$scratch = $ws.Range("Z1")

function Write-TextValue($range, [string]$text) {
    $escaped = $text.Replace('"', '""')
    $scratch.Formula = '="' + $escaped + '"'
    $scratch.Copy()
    $range.PasteSpecial(-4163)
}

# Append the new price-tracking row for 2026-02-07
Write-TextValue $ws.Range("A38") "2026-02-07"
Write-TextValue $ws.Range("B38") "699000"
Write-TextValue $ws.Range("C38") "38"
Write-TextValue $ws.Range("D38") "1"

$scratch.ClearContents()
